$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Rows("10").Delete()
$ws.Range("M28").Select()
